# TermBO fertig -> fertig für Tests schreiben
#
# Adds the new "Fragen" / "Bedingung" block (rows 77-86) to Tabelle1,
# matching the author's original edit order so the shared-string table
# is built up in the same sequence as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77 - intro sentence
$ws.Range("B77").Value = "wenn Sprache dargestellt wird, kann es sein, dass es entweder keine Specialty oder TechnicalTerm gibt (in einer Sprache), aber in einer anderen Sprache wurde die Beziehung schon definiert"

# Rows 83-86 - "Fragen:" block (written before the "Bedingung:" block below,
# matching the shared-string insertion order of the original commit)
$ws.Range("B83").Value = "Fragen:"
$ws.Range("C83").Value = "Warum gibt es laut den SQL Befehlen beim Aufbau der Tabellen und Beziehungen keine Verbindung zwischen Term und Translation?"
$ws.Range("C84").Value = "Befehle werden nicht weiterausgeführt bei TermDAOTest, hängt das mit detach zusammen bzw braucht es hierher ein Merge warum reicht keine neue Transaction?"
$ws.Range("C85").Value = "Wenn zB ein Term gelöscht wird, wird automatisch auch der Cascade Rest gelöscht (Translation)?"
$ws.Range("C86").Value = 'Fremdschlüsselbeziehung wird bei POJO als Objekt angegeben (zB bei Translation -> Language, Term) aber als ID des Objekts gespeichert: wie ist der Zugriff bei einer SQL Abfrage -> zB man will alle Translations einer Specialty'

# Rows 78-79 - "Bedingung:" block
$ws.Range("B78").Value = "Bedingung: "
$ws.Range("C78").Value = 'specialty darf niemals den gleichen "name" wie technicalTerm haben'
$ws.Range("C79").Value = "gleiche namen bei unterschiedlichen translation/sprachen sind möglich"

# Column E on the new rows keeps the centred column style (style index 1)
# but no content, mirroring the rest of the sheet's layout rows.
$ws.Range("E77").HorizontalAlignment = -4108
$ws.Range("E78").HorizontalAlignment = -4108
$ws.Range("E79").HorizontalAlignment = -4108
$ws.Range("E80").HorizontalAlignment = -4108
$ws.Range("E81").HorizontalAlignment = -4108
$ws.Range("E82").HorizontalAlignment = -4108
$ws.Range("E83").HorizontalAlignment = -4108
$ws.Range("E84").HorizontalAlignment = -4108
$ws.Range("E85").HorizontalAlignment = -4108
$ws.Range("E86").HorizontalAlignment = -4108

# Restore the selection/scroll position to where the author left off editing.
$ws.Range("D73").Select()
